$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 144; this shifts the old rows 144..226 down to 145..227,
# matching the rest of the diff (each former row's data now lives one row lower).
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new data record.
$ws.Cells.Item(144, 1).Value = 3
$ws.Cells.Item(144, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44960
$ws.Cells.Item(144, 5).Value = 5
$ws.Cells.Item(144, 6).Value = 100112052
$ws.Cells.Item(144, 7).Value = "Albahaca"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 155
$ws.Cells.Item(144, 11).Value = 4000
$ws.Cells.Item(144, 12).Value = 4500
$ws.Cells.Item(144, 13).Value = 4226
$ws.Cells.Item(144, 14).Value = "`$/docena de matas"
$ws.Cells.Item(144, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(144, 16).Value = 704
$ws.Cells.Item(144, 17).Value = 6
$ws.Cells.Item(144, 18).Value = "Hortaliza"
